# Auto-generated script to apply cryptos.xlsx diff via Excel COM-interop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.697.29"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.556.76"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.09"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.26"
$ws.Range("E6").Value = "  +7.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.49"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  +9.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.52"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "2.494.29"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.880"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.49"
$ws.Range("E16").Value = "  +3.38%  "
$ws.Range("D17").Value = "42.756.13"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.60"
$ws.Range("E18").Value = "  +9.38%  "
$ws.Range("D19").Value = "0.0₃0987"
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.59"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.55"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "255.87"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "28.03"
$ws.Range("E25").Value = "  -4.02%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.07"
$ws.Range("E27").Value = "  +9.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.10"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.99"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.75"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0801"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.33"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("B36").Value = "EnergySwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.58"
$ws.Range("E36").Value = "  +9.49%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.27"
$ws.Range("E37").Value = "  +17.29%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.05"
$ws.Range("E41").Value = "  +29.70%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.36"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0305"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D45").Value = "2.058.65"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.05"
$ws.Range("E46").Value = "  +4.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.23"
$ws.Range("E47").Value = "  +6.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.79"
$ws.Range("E48").Value = "  +11.74%  "
$ws.Range("D49").Value = "2.809.57"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.75"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("E51").Value = "  +3.27%  "
